$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Visual Studio 2017" row (row 5), matching the order the
# author originally typed things in (tool name, then the URL/hyperlink
# column, then the tool type, then the three long-text columns).
$ws.Range("A5").Value = "Visual Studio 2017"
$ws.Range("C5").Value = "https://visualstudio.microsoft.com/downloads/"
$ws.Range("B5").Value = "Debugger/IDE"
$ws.Range("D5").Value = "Using the debugger with an executable can take a bit of time to get used to, but once you understand the basics it becomes easier."
$ws.Range("E5").Value = "This program allows the user to step through (trace) each line of the given  assembly code of the executable.  This can be very helpful to see how inputting different values changes the flow of the program."
$ws.Range("F5").Value = "Knowing how to upload an executable so that you can trace through it is not obvious and takes some time to find."

# Turn the URL in C5 into a real hyperlink, same as the other tools above it.
$ws.Hyperlinks.Add($ws.Range("C5"), "https://visualstudio.microsoft.com/downloads/")

# The three "long text" cells should wrap, like the corresponding cells in
# the rows above.
$ws.Range("D5:F5").WrapText = $true

# The row grew tall enough to show all of the wrapped text.
$ws.Rows.Item(5).RowHeight = 111.6

# Column C (URLs) was widened/auto-fit to comfortably show the new, longer
# Visual Studio download link.
$ws.Columns.Item(3).AutoFit()

# Reflect where the user ended up looking after adding the row.
$ws.Range("F5").Select() | Out-Null
